$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value2 = '27.374.37'
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Value2 = '  -0.27%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value2 = '1.717.70'
$ws.Cells.Item(3, 2).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)
$ws.Cells.Item(3, 5).Value2 = '  -0.29%  '

$ws.Cells.Item(4, 5).Value2 = '  +0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value2 = '224.84'
$ws.Cells.Item(5, 2).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)
$ws.Cells.Item(5, 5).Value2 = '  -0.20%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value2 = '0.5308'
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4122)
$ws.Cells.Item(6, 5).Value2 = '  -0.56%  '

$ws.Cells.Item(7, 5).Value2 = '  -0.02%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value2 = '0.06717'
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)
$ws.Cells.Item(8, 5).Value2 = '  +2.30%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value2 = '0.2661'
$ws.Cells.Item(9, 2).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4122)
$ws.Cells.Item(9, 5).Value2 = '  +0.24%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value2 = '20.94'
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4122)
$ws.Cells.Item(10, 5).Value2 = '  -2.51%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value2 = '0.07688'
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4122)
$ws.Cells.Item(11, 5).Value2 = '  +0.37%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value2 = '4.488'
$ws.Cells.Item(12, 2).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4122)
$ws.Cells.Item(12, 5).Value2 = '  -2.09%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value2 = '1.953.38'
$ws.Cells.Item(13, 2).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4122)
$ws.Cells.Item(13, 5).Value2 = '  -0.31%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value2 = '1.713.48'
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4122)
$ws.Cells.Item(14, 5).Value2 = '  -0.61%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value2 = '0.5824'
$ws.Cells.Item(15, 2).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4122)
$ws.Cells.Item(15, 5).Value2 = '  +0.83%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value2 = '0.0₅8222'
$ws.Cells.Item(16, 2).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4122)
$ws.Cells.Item(16, 5).Value2 = '  -0.47%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value2 = '68.05'
$ws.Cells.Item(17, 2).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4122)
$ws.Cells.Item(17, 5).Value2 = '  +0.50%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value2 = '27.366.96'
$ws.Cells.Item(18, 2).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4122)
$ws.Cells.Item(18, 5).Value2 = '  -0.33%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value2 = '223.83'
$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4122)
$ws.Cells.Item(19, 5).Value2 = '  +2.62%  '

$ws.Cells.Item(20, 5).Value2 = '  +0.00%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value2 = '4.668'
$ws.Cells.Item(21, 2).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4122)
$ws.Cells.Item(21, 5).Value2 = '  -0.96%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value2 = '10.47'
$ws.Cells.Item(22, 2).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4122)
$ws.Cells.Item(22, 5).Value2 = '  -0.78%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value2 = '6.033'
$ws.Cells.Item(23, 2).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4122)
$ws.Cells.Item(23, 5).Value2 = '  +0.21%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value2 = '1.005'
$ws.Cells.Item(24, 2).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4122)
$ws.Cells.Item(24, 5).Value2 = '  -0.09%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value2 = '144.41'
$ws.Cells.Item(25, 2).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4122)
$ws.Cells.Item(25, 5).Value2 = '  +1.01%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value2 = '1.702'
$ws.Cells.Item(26, 2).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4122)
$ws.Cells.Item(26, 5).Value2 = '  -2.74%  '

$ws.Cells.Item(27, 5).Value2 = '  -1.74%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value2 = '7.250'
$ws.Cells.Item(28, 2).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4122)
$ws.Cells.Item(28, 5).Value2 = '  -0.95%  '

$ws.Cells.Item(29, 5).Value2 = '  -1.00%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value2 = '0.05401'
$ws.Cells.Item(30, 2).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4122)
$ws.Cells.Item(30, 5).Value2 = '  -1.34%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value2 = '1.295'
$ws.Cells.Item(31, 2).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4122)
$ws.Cells.Item(31, 5).Value2 = '  -0.31%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value2 = '3.487'
$ws.Cells.Item(32, 2).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4122)
$ws.Cells.Item(32, 5).Value2 = '  -1.78%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value2 = '3.422'
$ws.Cells.Item(33, 2).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4122)
$ws.Cells.Item(33, 5).Value2 = '  -0.21%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value2 = '1.635'
$ws.Cells.Item(34, 2).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4122)
$ws.Cells.Item(34, 5).Value2 = '  -1.14%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value2 = '2.861'
$ws.Cells.Item(35, 2).Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4122)
$ws.Cells.Item(35, 5).Value2 = '  +0.21%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value2 = '0.9563'
$ws.Cells.Item(36, 2).Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4122)
$ws.Cells.Item(36, 5).Value2 = '  +0.17%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value2 = '2.395'
$ws.Cells.Item(37, 2).Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4122)
$ws.Cells.Item(37, 5).Value2 = '  -1.20%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value2 = '0.5905'
$ws.Cells.Item(38, 2).Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4122)
$ws.Cells.Item(38, 5).Value2 = '  -0.37%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value2 = '1.147.37'
$ws.Cells.Item(39, 2).Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4122)
$ws.Cells.Item(39, 5).Value2 = '  +9.56%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value2 = '0.01654'
$ws.Cells.Item(40, 2).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4122)
$ws.Cells.Item(40, 5).Value2 = '  +0.61%  '

$ws.Cells.Item(41, 5).Value2 = '  -1.49%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value2 = '1.005'
$ws.Cells.Item(42, 2).Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4122)
$ws.Cells.Item(42, 5).Value2 = '  +0.02%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value2 = '0.8421'
$ws.Cells.Item(43, 2).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4122)
$ws.Cells.Item(43, 5).Value2 = '  -0.71%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value2 = '101.05'
$ws.Cells.Item(44, 2).Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4122)
$ws.Cells.Item(44, 5).Value2 = '  -0.16%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value2 = '1.860.40'
$ws.Cells.Item(45, 2).Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4122)
$ws.Cells.Item(45, 5).Value2 = '  -0.36%  '

$ws.Cells.Item(46, 5).Value2 = '  -6.83%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value2 = '57.92'
$ws.Cells.Item(47, 2).Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4122)
$ws.Cells.Item(47, 5).Value2 = '  -1.17%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value2 = '0.4587'
$ws.Cells.Item(48, 2).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4122)
$ws.Cells.Item(48, 5).Value2 = '  +2.36%  '

$ws.Cells.Item(49, 2).Value2 = 'Frax'
$ws.Cells.Item(49, 3).Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value2 = '1.005'
$ws.Cells.Item(49, 2).Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4122)
$ws.Cells.Item(49, 5).Value2 = '  +0.27%  '

$ws.Cells.Item(50, 2).Value2 = 'EnergySwap'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value2 = '8.134'
$ws.Cells.Item(50, 2).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4122)
$ws.Cells.Item(50, 5).Value2 = '  -0.38%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value2 = '0.05202'
$ws.Cells.Item(51, 2).Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4122)
$ws.Cells.Item(51, 5).Value2 = '  -0.87%  '
